$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K6").Value = 31
$ws.Range("L6").Value = 81.65000000000001

$ws.Range("K7").Value = 135.838
$ws.Range("L7").Value = 102.54

$ws.Range("K8").Value = 208.948
$ws.Range("L8").Value = 251.47
